# Reproduce the edits captured in the target diff:
#  - Sheet "5.3": C7/C8 become formulas (=724 / =788) instead of static 7.24/7.88
#  - Sheet "5.3": zoom to 190%, selection moves to B2, no longer the active tab
#  - Sheet "5.1": becomes the active tab, selection moves to D13
#  - Workbook window: resize/reposition (maximized-looking geometry)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("5.1")
$ws3 = $wb.Worksheets.Item("5.3")

# --- Data edit on "5.3": C7 -> =724, C8 -> =788 ---------------------------
$ws3.Range("C7").Formula = "=724"
$ws3.Range("C8").Formula = "=788"

# --- View state on "5.3": zoom 190%, select B2 ----------------------------
$ws3.Activate()
$ws3.Range("B2").Select()
$excel.ActiveWindow.Zoom = 190
$excel.ActiveWindow.ZoomScaleNormal = 190

# --- View state on "5.1": becomes the active/selected sheet, select D13 --
$ws1.Activate()
$ws1.Range("D13").Select()

# --- Workbook window geometry (maximized-style position/size) ------------
$win = $excel.ActiveWindow
$win.WindowState = -4137   # xlMaximized
$win.Left   = -108
$win.Top    = -108
$win.Width  = 23256
$win.Height = 12576
